$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying model output table stores its numeric-looking figures as
# plain text (General-formatted cells containing shared-string text, not
# real numbers). Assigning a numeric-looking string straight to .Value
# would make Excel silently re-interpret it as a Number, which changes the
# cell's XML type (t="s" -> no t / numeric <v>) and is NOT what the target
# edit does. To keep these cells as text we temporarily borrow a
# text-formatted donor cell's number format (copy/paste formats), set the
# value (now safely kept as text), then paste back the formatting of a
# cell that already has the right, original look - this avoids leaving
# any new/unused style behind.

function Set-TextValue($targetAddr, $value, $styleDonorAddr) {
    # Step 1: temporarily apply a text-capable (@) number format so the
    # value assigned below is retained as text rather than coerced to a
    # number.
    $ws.Range("F8").Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    # Step 2: write the literal text.
    $ws.Range($targetAddr).Value = $value

    # Step 3: restore the exact original formatting/style of this cell.
    $ws.Range($styleDonorAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# --- Row 2: intrcpt ---
Set-TextValue "D2" "0.479" "A2"
Set-TextValue "F2" "0.6317" "F8"

# --- Row 3: Diet_HCOherbivore ---
Set-TextValue "B3" " 0.012" "A2"
Set-TextValue "D3" "3.297" "A2"
Set-TextValue "F3" "0.1923" "F8"

# --- Row 4: Diet_HCOomnivore ---
Set-TextValue "B4" "-0.124" "A2"
Set-TextValue "C4" "0.077" "A2"
Set-TextValue "D4" "3.297" "A2"
Set-TextValue "F4" "0.1923" "F8"

# --- Row 5: Migratmigrant ---
Set-TextValue "C5" "0.063" "A2"
Set-TextValue "D5" "0.927" "A2"
Set-TextValue "F5" "0.3357" "F8"

# --- Row 6: GenLength_y_IUCN.y ---
Set-TextValue "D6" "0.989" "A2"
Set-TextValue "F6" "0.3200" "F8"

# --- Row 7: abs_lat ---
Set-TextValue "D7" "0.073" "A2"
Set-TextValue "F7" "0.7873" "F8"

# --- Remove row 8 (the old standalone "Pvalue" covariate row); CG model
# now folds PdeltaAIC in as a covariate instead of a separate row. ---
$ws.Rows.Item(8).Delete()
